# Applies the content edits described by the commit diff to 论文记录.xlsx
# (standardizing many "Data" column dataset-name strings, fixing one paper
# title typo, and clearing a couple of now-superfluous cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Cell value edits -------------------------------------------------

$ws.Range("J3").Value = 'WSJ'
$ws.Range("J9").Value = 'Yapex, GENIA, UMLS, Mammalian Phenotype Ontology'
$ws.Range("J10").Value = 'Others'
$ws.Range("J11").Value = 'CoNLL 2003'
$ws.Range("A12").Value = 'Exploiting Wikipedia as External Knowledge for Named Entity Recognition'
$ws.Range("J12").Value = 'CoNLL 2003, Wikipedia '
$ws.Range("J13").Value = 'Others'
$ws.Range("J15").Value = 'CoNLL 2003, MUC7, Others'
$ws.Range("J18").Value = 'Penn Treebank, CoNLL 2000, WSJ'
$ws.Range("J21").Value = 'MeSH,OMIM'
$ws.Range("J22").Value = 'DBpedia'
$ws.Range("J24").Value = 'Others'
$ws.Range("J26").Value = 'CoNLL 2003, CoNLL 2000 '
$ws.Range("J29").Value = 'Wikipedia， CoNLL 2002, CoNLL 2003, Others'
$ws.Range("J30").Value = 'CoNLL 2003, ACE 2, ACE 2003, MUC7, Wikipedia, Gigaword5'
$ws.Range("J31").Value = 'Tweets, ACE 2005, Enron,CoNLL 2003 '
$ws.Range("J32").Value = 'SCAI'
$ws.Range("J33").Value = 'Wikipedia, Web Pages, Tweets'
$ws.Range("J35").Value = 'Sina Weibo'
$ws.Range("J36").Value = 'Web Pages'
$ws.Range("J37").Value = 'Tweets, CoNLL 2003 '
$ws.Range("J38").Value = $null
$ws.Range("J39").Value = 'Tweets'
$ws.Range("J40").Value = 'CoNLL 2003'
$ws.Range("J41").Value = 'CoNLL 2003, Others'
$ws.Range("J42").Value = 'CoNLL 2002, CoNLL 2003'
$ws.Range("J44").Value = 'CoNLL 2003'
$ws.Range("J45").Value = 'CoNLL 2002, CoNLL 2003 '
$ws.Range("J46").Value = 'Tweets'
$ws.Range("J47").Value = 'WordNet, CoNLL 2003 '
$ws.Range("J48").Value = "TAC KBP 2013, CoNLL`n2003, CoNLL 2002, Ontonotes,Wikipedia"
$ws.Range("J50").Value = "WordSim`n353, Others"
$ws.Range("J51").Value = 'WSJ, CoNLL 2003, Wikipedia'
$ws.Range("J52").Value = 'SIGHAN 2005 shared task, Sina Weibo'
$ws.Range("J53").Value = 'CoNLL 2003'
$ws.Range("J55").Value = 'Wikipedia, Others'
$ws.Range("J56").Value = 'North American News'
$ws.Range("J57").Value = 'Others'

# --- Cosmetic view-state update (scroll position / active cell) ------

$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Application.ActiveWindow.ScrollRow = 52
$ws.Range("J56").Select()
